$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the schedule validation note to also mention delete/update, and mark it checked
$ws.Range("A29").Value = "validasi admin doang yg bisa add/delete/update skedul"
$ws.Range("B29").Value = "check"

# Mark the "bikin form change password" row as checked (change password form added)
$ws.Range("B45").Value = "check"

# Move the selection/view to B10 (also clears the old scrolled topLeftCell state)
$ws.Range("B10").Select()
